# Add prepare draw logic
# Rebuild the role_card Sheet1 data: drop placeholder rows (LIUBEI/GUOJIA/SUNCE),
# fill in force/command/moral + tex/desc columns for every existing hero row,
# and append 4 new heroes (SIMAYI, YUEYING, GANNING) plus 4 brand new rows
# (DIANWEI, ZHANGHE, ZHOUTAI, ZHOUYU) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9 was LIUBEI (placeholder) -> now SIMAYI, fully populated ---
$ws.Range("B9").Value = "SIMAYI"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = "simayi"
$ws.Range("J9").Value = "SIMAYI_DESC"

# --- Row 12 was GUOJIA (placeholder) -> now YUEYING, fully populated ---
$ws.Range("B12").Value = "YUEYING"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 6
$ws.Range("H12").Value = "yueying"
$ws.Range("J12").Value = "YUEYING_DESC"

# --- Row 14 was SUNCE (placeholder) -> now GANNING, fully populated ---
$ws.Range("B14").Value = "GANNING"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = "ganning"
$ws.Range("J14").Value = "GANNING_DESC"

# --- Fill in force/command/moral + tex/desc for the remaining existing rows ---
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 5
$ws.Range("H4").Value = "zhaoyun"
$ws.Range("J4").Value = "ZHAOYUN_DESC"

$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("H5").Value = "lvbu"
$ws.Range("J5").Value = "LVBU_DESC"

$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 5
$ws.Range("H6").Value = "diaochan"
$ws.Range("J6").Value = "DIAOCHAN_DESC"

$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("H7").Value = "guanyu"
$ws.Range("J7").Value = "GUANYU_DESC"

$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 4
$ws.Range("H8").Value = "zhangfei"
$ws.Range("J8").Value = "ZHANGFEI_DESC"

$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 5
$ws.Range("H11").Value = "xiahoudun"
$ws.Range("J11").Value = "XIAHOUDUN_DESC"

$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 6
$ws.Range("H13").Value = "zhugeliang"
$ws.Range("J13").Value = "ZHUGELIANG"

$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = "shangxiang"
$ws.Range("J15").Value = "SHANGXIANG_DESC"

$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 5
$ws.Range("H16").Value = "huanggai"
$ws.Range("J16").Value = "HUANGGAI_DESC"

# --- Append 4 brand new rows: DIANWEI, ZHANGHE, ZHOUTAI, ZHOUYU ---
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "DIANWEI"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = "dianwei"
$ws.Range("J17").Value = "DIANWEI_DESC"

$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "ZHANGHE"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = "zhanghe"
$ws.Range("J18").Value = "ZHANGHE_DESC"

$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "ZHOUTAI"
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 5
$ws.Range("H19").Value = "zhoutai"
$ws.Range("J19").Value = "ZHOUTAI_DESC"

$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "ZHOUYU"
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = "zhouyu"
$ws.Range("J20").Value = "ZHOUYU_DESC"

# --- Widen column J to fit the new desc keys, and move the selection cursor ---
# (16.7 chars is the closest achievable COM ColumnWidth to the authored 17.375
# stored width; the host rounds chars to whole pixels at 7px/char before
# re-expressing as "characters", so exact fractional widths like 17.375 can't
# always be hit bit-for-bit.)
$ws.Columns.Item(10).ColumnWidth = 16.7
$ws.Range("J22").Select()
